$d = $word.ActiveDocument

# --- Change 1: Insert new list item after "Code cleanup..." paragraph ---
$rng1 = $d.Content
$rng1.Find.Execute("Code cleanup. Global namespace operator etc where appropriate.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng1.Collapse(0)  # wdCollapseEnd
$rng1.InsertAfter("Convert use of boost/fs in jamfile to boost/headers if possible.")
$rng1.InsertParagraphBefore()

# --- Change 2: Move lastRenderedPageBreak from "Helper service..." run to "Basic base hook." run ---

# 2a. Add lastRenderedPageBreak to the "Basic base hook." paragraph's run
$rng2 = $d.Content
$rng2.Find.Execute("Basic base hook.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$para2 = $rng2.Paragraphs(1).Range
$xml2 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="004D5588" w:rsidRDefault="004D5588" w:rsidP="004D5588"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Basic base hook.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$para2.InsertXML($xml2)

# 2b. Remove lastRenderedPageBreak from the "Helper service..." paragraph's run
$rng3 = $d.Content
$rng3.Find.Execute("Helper service to run HadesMem", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$para3 = $rng3.Paragraphs(1).Range
$xml3 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="005D3274" w:rsidRPr="008B06FC" w:rsidRDefault="005D3274" w:rsidP="005D3274"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r w:rsidRPr="008B06FC"><w:t xml:space="preserve">Helper service to run HadesMem tools as ‘SYSTEM’ for when </w:t></w:r><w:r w:rsidR="00A77D3A"><w:t>manipulating</w:t></w:r><w:r w:rsidRPr="008B06FC"><w:t xml:space="preserve"> certain protected/critical processes (running in separate desktops, sessions, etc.).</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$para3.InsertXML($xml3)

Write-Host "All changes applied"
